$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '23.374.61'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.83%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.626.05'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -1.00%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9997'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.9998'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.09%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '304.53'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -1.31%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3785'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '51.89'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -2.24%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3628'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -1.74%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.229'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -4.31%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -1.33%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.9997'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.12%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.69'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -2.56%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.551'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -1.98%  '
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -3.43%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.222'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -3.62%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.632.55'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.55%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '93.55'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.60%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06902'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.81%  '
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -3.38%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9998'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.14%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.416'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -2.86%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '23.382.91'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.77%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.72'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -2.31%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.245'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +3.10%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +1.13%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '149.89'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -1.22%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.289'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -1.05%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '134.08'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -1.68%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.304'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -5.03%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.812.67'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.25%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.780'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -1.22%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '10.99'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +5.12%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9520'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02784'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -1.46%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.2522'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.08820'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.57%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.098'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -2.38%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.07128'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -4.98%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.359'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -3.21%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.7066'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -1.82%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '16.18'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '12.28'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -3.24%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6449'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -3.15%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.320'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -2.16%  '
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.10%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.993'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -1.38%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.07990'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -1.09%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -1.60%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '125.72'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -4.47%  '
